# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" summary row at the top of the "总计" sheet
#    (pushing the existing "2021-Q4" summary row down to row 3).
# 2. Duplicate the "2021-Q4" detail sheet (to inherit its header/row
#    formatting and page setup), place the copy right before "2021-Q4",
#    rename it "2022-Q4", and replace its contents with the 2022-Q4
#    fund-holding data (adding 2 extra data rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert new row 2 for 2022-Q4, shift 2021-Q4 down
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("A2:D2").Insert()

# carry the bordered/bold "A-column" style down onto the new row
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.03

# the row-index column (A) renumbers: the old 2021-Q4 row becomes index 1
$summary.Range("A3").Value = 1

# ---------------------------------------------------------------------
# Step 2: create the "2022-Q4" detail sheet (copy of "2021-Q4" so the
# header/row styling and sheet page setup match), positioned before it.
# ---------------------------------------------------------------------
$detail2021 = $wb.Worksheets.Item("2021-Q4")
$detail2021.Copy($detail2021)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q4"

# extend the styled A-column down to rows 5 and 6 (only 4 rows existed)
$newSheet.Range("A4").Copy()
$newSheet.Range("A5:A6").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'007277"
$newSheet.Range("C2").Value = "恒生前海消费升级混合"
$newSheet.Range("D2").Value = "'0.47"
$newSheet.Range("E2").Value = "'85.02"
$newSheet.Range("F2").Value = "'3.31"
$newSheet.Range("G2").Value = "'0.0156"
$newSheet.Range("H2").Value = 6
$newSheet.Range("B2:G2").ClearFormats()

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'008861"
$newSheet.Range("C3").Value = "西部利得港股通新机遇灵活配置混合A"
$newSheet.Range("D3").Value = "'0.25"
$newSheet.Range("E3").Value = "'87.69"
$newSheet.Range("F3").Value = "'3.38"
$newSheet.Range("G3").Value = "'0.0084"
$newSheet.Range("H3").Value = 10
$newSheet.Range("B3:G3").ClearFormats()

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'010093"
$newSheet.Range("C4").Value = "西部利得港股通新机遇灵活配置混合C"
$newSheet.Range("D4").Value = "'0.12"
$newSheet.Range("E4").Value = "'87.69"
$newSheet.Range("F4").Value = "'3.38"
$newSheet.Range("G4").Value = "'0.0041"
$newSheet.Range("H4").Value = 10
$newSheet.Range("B4:G4").ClearFormats()

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'013182"
$newSheet.Range("C5").Value = "安信港股通精选混合C"
$newSheet.Range("D5").Value = "'0.12"
$newSheet.Range("E5").Value = "'69.28"
$newSheet.Range("F5").Value = "'2.54"
$newSheet.Range("G5").Value = "'0.0030"
$newSheet.Range("H5").Value = 9
$newSheet.Range("B5:G5").ClearFormats()

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'013181"
$newSheet.Range("C6").Value = "安信港股通精选混合A"
$newSheet.Range("D6").Value = "'0.02"
$newSheet.Range("E6").Value = "'69.28"
$newSheet.Range("F6").Value = "'2.54"
$newSheet.Range("G6").Value = "'0.0005"
$newSheet.Range("H6").Value = 9
$newSheet.Range("B6:G6").ClearFormats()

# restore the originally-active sheet/tab ("总计")
$summary.Activate()
